$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Bat', ['Token Creature — Bat', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Cleric', ['Token Creature — Cleric', '{3}{W}{B}{B}, {T}, Sacrifice this creature: Return a card named Deathpact Angel from your graveyard to the battlefield.', '1/1'])"
$ws.Range("A5").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '6/6'])"
$ws.Range("A6").Value = "('Goblin', ['Token Creature — Goblin', '2/1'])"
$ws.Range("A7").Value = "('Ooze', ['Token Creature — Ooze', '*/*'])"
$ws.Range("A8").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"
$ws.Range("A9").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A10").Value = "('Wurm', ['Token Creature — Wurm', '6/6'])"

$ws.Range("A11:A33").EntireRow.Delete()
